$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'273.87"
$ws.Range('E2').Value = "'-1.65%"
$ws.Range('D3').Value = "'26.65"
$ws.Range('E3').Value = "'-2.40%"
$ws.Range('D4').Value = "'4.871"
$ws.Range('E4').Value = "'1.15%"
$ws.Range('D5').Value = "'0.06320"
$ws.Range('E5').Value = "'0.58%"
$ws.Range('D6').Value = "'6.880"
$ws.Range('E6').Value = "'0.35%"
$ws.Range('B7').Value = "'GateToken"
$ws.Range('C7').Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range('D7').Value = "'3.332"
$ws.Range('E7').Value = "'1.86%"
$ws.Range('B8').Value = "'FTXToken"
$ws.Range('C8').Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range('D8').Value = "'1.216"
$ws.Range('E8').Value = "'27.64%"
$ws.Range('B9').Value = "'MXToken"
$ws.Range('C9').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D9').Value = "'0.8701"
$ws.Range('E9').Value = "'-0.79%"
$ws.Range('B10').Value = "'WazirX"
$ws.Range('C10').Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range('D10').Value = "'0.1455"
$ws.Range('E10').Value = "'-0.03%"
$ws.Range('B11').Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range('C11').Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range('D11').Value = "'0.05080"
$ws.Range('E11').Value = "'-1.64%"
$ws.Range('B12').Value = "'MandalaExchangeToken"
$ws.Range('C12').Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range('D12').Value = "'0.07382"
$ws.Range('E12').Value = "'1.30%"
$ws.Range('B13').Value = "'BitrueCoin"
$ws.Range('C13').Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range('D13').Value = "'0.02969"
$ws.Range('E13').Value = "'-6.26%"
$ws.Range('B14').Value = "'BitMartToken"
$ws.Range('C14').Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range('D14').Value = "'0.09037"
$ws.Range('E14').Value = "'-0.19%"
$ws.Range('B15').Value = "'BitForexToken"
$ws.Range('C15').Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range('D15').Value = "'0.001572"
$ws.Range('E15').Value = "'1.25%"
$ws.Range('B16').Value = "'One"
$ws.Range('C16').Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range('D16').Value = "'0.0006279"
$ws.Range('E16').Value = "'0.36%"
$ws.Range('B17').Value = "'TigerCash"
$ws.Range('C17').Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range('D17').Value = "'0.006023"
$ws.Range('E17').Value = "'1.59%"
$ws.Range('B18').Value = "'LEO"
$ws.Range('C18').Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range('D18').Value = "'3.453"
$ws.Range('E18').Value = "'-0.43%"
$ws.Range('D19').Value = "'2.284"
$ws.Range('E19').Value = "'0.79%"
$ws.Range('E20').Value = "'2.55%"
$ws.Range('E21').Value = "'1.25%"
$ws.Range('D22').Value = "'3.902"
$ws.Range('E22').Value = "'1.39%"
$ws.Range('D23').Value = "'0.04387"
$ws.Range('E23').Value = "'1.78%"
$ws.Range('D24').Value = "'0.001174"
$ws.Range('E24').Value = "'-0.03%"
$ws.Range('D25').Value = "'0.004267"
$ws.Range('E25').Value = "'-0.24%"
$ws.Range('D26').Value = "'0.0001199"
$ws.Range('E26').Value = "'0.04%"
$ws.Range('D27').Value = "'0.0001693"
$ws.Range('E27').Value = "'-4.59%"
$ws.Range('D40').Value = "'0.04030"
$ws.Range('E40').Value = "'-0.03%"
$ws.Range('D41').Value = "'0.006727"
$ws.Range('E41').Value = "'0.22%"
$ws.Range('D42').Value = "'0.1167"
$ws.Range('D43').Value = "'0.002099"
$ws.Range('E43').Value = "'-0.10%"
$ws.Range('D44').Value = "'0.01254"
$ws.Range('E44').Value = "'-10.89%"
$ws.Range('D45').Value = "'0.00005305"
$ws.Range('E45').Value = "'2.36%"
$ws.Range('B46').Value = "'BOLO"
$ws.Range('C46').Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range('D46').Value = "'2.357"
$ws.Range('E46').Value = "'2.03%"
$ws.Range('B47').Value = "'CoinbaseStockToken"
$ws.Range('C47').Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range('D47').Value = "'0.02000"
$ws.Range('E47').Value = "'-33.06%"
